$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 9 (pushes GPU_small_e_czech etc. down by one row)
$ws.Rows.Item(9).Insert()

# Clear any inherited formatting on the new row's table cells so we can set it explicitly
$ws.Range("D9:H9").ClearFormats()

# Set the new model name in column D
$ws.Range("D9").Value2 = "GPU_bert_cased"

$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138
$xlEdgeLeft = 7
$xlEdgeRight = 10

# Re-create the table border look for the new row: thin border around every cell ...
$ws.Range("D9:H9").Borders.LineStyle = $xlContinuous
$ws.Range("D9:H9").Borders.Weight = $xlThin

# ... and thicken the two outer edges of the table (left of column D, right of column H)
$ws.Range("D9").Borders.Item($xlEdgeLeft).Weight = $xlMedium
$ws.Range("H9").Borders.Item($xlEdgeRight).Weight = $xlMedium

# Leave the cursor where the author left it
$ws.Range("K9").Select()
